# The "States" worksheet had two columns removed from the data table:
# "Thousands Of Residents" (col D) and "Household Median Income" (col E).
# The two columns that used to be F ("Sales Per Thousand Residents") and
# G ("Median Sales Price Vs Median Household Income") slide left to
# become the new D and E columns, and the two now-unused strings are
# dropped from the shared string table automatically by Excel since no
# cell references them any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns D and E entirely (shifts F,G left into D,E).
$ws.Columns("D:E").Delete()
